$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 1755.75
$ws.Range("I38").Value = 1341
$ws.Range("K38").Value = 4023
$ws.Range("M38").Value = -3651
$ws.Range("H86").Value = 1469.3846
$ws.Range("I86").Value = 1149.375
$ws.Range("J86").Value = 1981.4
$ws.Range("K86").Value = 1149.375
$ws.Range("L86").Value = 1981.4
$ws.Range("M86").Value = -26.375
$ws.Range("N86").Value = -4227.4
$ws.Range("H89").Value = 1469.3846
$ws.Range("I89").Value = 1149.375
$ws.Range("J89").Value = 1981.4
$ws.Range("K89").Value = 5746.875
$ws.Range("L89").Value = 9907
$ws.Range("M89").Value = -130.875
$ws.Range("N89").Value = -21139
$ws.Range("H112").Value = 3532
$ws.Range("J112").Value = 3532
$ws.Range("L112").Value = 10596
$ws.Range("N112").Value = -12812
$ws.Range("H138").Value = 2320.8333
$ws.Range("J138").Value = 2282.5
$ws.Range("L138").Value = 6847.5
$ws.Range("N138").Value = -17127.5

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3652.8
$ws.Range("I32").Value = 2153.606
$ws.Range("J32").Value = 7775.5835
$ws.Range("K32").Value = 2153.606
$ws.Range("L32").Value = 7775.5835
$ws.Range("M32").Value = -1866.606
$ws.Range("N32").Value = -8349.583500000001
$ws.Range("H74").Value = 4640.1816
$ws.Range("I74").Value = 4682.3335
$ws.Range("J74").Value = 4450.5
$ws.Range("K74").Value = 4682.3335
$ws.Range("L74").Value = 4450.5
$ws.Range("M74").Value = -3808.3335
$ws.Range("N74").Value = -6198.5
$ws.Range("H77").Value = 4640.1816
$ws.Range("I77").Value = 4682.3335
$ws.Range("J77").Value = 4450.5
$ws.Range("K77").Value = 23411.6675
$ws.Range("L77").Value = 22252.5
$ws.Range("M77").Value = -19043.6675
$ws.Range("N77").Value = -30988.5
$ws.Range("H132").Value = 2208.074
$ws.Range("I132").Value = 1351
$ws.Range("J132").Value = 3279.4167
$ws.Range("K132").Value = 4053
$ws.Range("L132").Value = 9838.250100000001
$ws.Range("M132").Value = -1523
$ws.Range("N132").Value = -14898.2501

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 605
$ws.Range("I64").Value = 594.7143
$ws.Range("J64").Value = 619.4
$ws.Range("K64").Value = 594.7143
$ws.Range("L64").Value = 619.4
$ws.Range("M64").Value = -369.7143
$ws.Range("N64").Value = -1069.4
$ws.Range("H67").Value = 605
$ws.Range("I67").Value = 594.7143
$ws.Range("J67").Value = 619.4
$ws.Range("K67").Value = 594.7143
$ws.Range("L67").Value = 619.4
$ws.Range("M67").Value = 185.2857
$ws.Range("N67").Value = -2179.4
$ws.Range("H107").Value = 2591.0527
$ws.Range("J107").Value = 3608.3333
$ws.Range("L107").Value = 3608.3333
$ws.Range("N107").Value = -7448.3333
$ws.Range("H126").Value = 40000
$ws.Range("J126").Value = 40000
$ws.Range("L126").Value = 40000
$ws.Range("N126").Value = -49880

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 129.16667
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()
$ws.Range("H31").Value = 2148.5806
$ws.Range("I31").Value = 983.913
$ws.Range("K31").Value = 983.913
$ws.Range("M31").Value = -688.913
$ws.Range("H34").Value = 2148.5806
$ws.Range("I34").Value = 983.913
$ws.Range("K34").Value = 983.913
$ws.Range("M34").Value = -781.913
$ws.Range("H59").Value = 39600
$ws.Range("J59").Value = 39600
$ws.Range("L59").Value = 39600
$ws.Range("N59").Value = -41890
$ws.Range("H99").Value = 2049
$ws.Range("I99").Value = 1433
$ws.Range("K99").Value = 1433
$ws.Range("M99").Value = 65
$ws.Range("H126").Value = 2049
$ws.Range("I126").Value = 1433
$ws.Range("K126").Value = 4299
$ws.Range("M126").Value = -1829
$ws.Range("H134").Value = 1145.8572
$ws.Range("I134").Value = 1145.8572
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 3437.5716
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -902.5715999999998
$ws.Range("N134").ClearContents()
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").ClearContents()
$ws.Range("N140").Value = 0

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 274
$ws.Range("J26").Value = 248
$ws.Range("L26").Value = 744
$ws.Range("N26").Value = -1320
$ws.Range("H55").Value = 52502
$ws.Range("J55").Value = 5000
$ws.Range("L55").Value = 15000
$ws.Range("N55").Value = -15354
$ws.Range("H131").Value = 7473853
$ws.Range("J131").Value = 12405.85
$ws.Range("L131").Value = 37217.55
$ws.Range("N131").Value = -47297.55
$ws.Range("H132").Value = 2085.7144
$ws.Range("I132").Value = 1600
$ws.Range("K132").Value = 14400
$ws.Range("M132").Value = -11870

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H55").Value = 9000
$ws.Range("J55").Value = 9000
$ws.Range("L55").Value = 9000
$ws.Range("N55").Value = -9654
$ws.Range("H70").Value = 4774.25
$ws.Range("I70").Value = 4983
$ws.Range("J70").Value = 4148
$ws.Range("K70").Value = 4983
$ws.Range("L70").Value = 4148
$ws.Range("M70").Value = -4713
$ws.Range("N70").Value = -4688
$ws.Range("H73").Value = 4774.25
$ws.Range("I73").Value = 4983
$ws.Range("J73").Value = 4148
$ws.Range("K73").Value = 4983
$ws.Range("L73").Value = 4148
$ws.Range("M73").Value = -4047
$ws.Range("N73").Value = -6020
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").ClearContents()
$ws.Range("N86").Value = 0
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").ClearContents()
$ws.Range("N89").Value = 0
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").ClearContents()
$ws.Range("N124").Value = 0

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5487.76
$ws.Range("I40").Value = 2807.3076
$ws.Range("J40").Value = 8391.583000000001
$ws.Range("K40").Value = 2807.3076
$ws.Range("L40").Value = 8391.583000000001
$ws.Range("M40").Value = -2671.3076
$ws.Range("N40").Value = -8663.583000000001
$ws.Range("H122").Value = 4248.387
$ws.Range("I122").Value = 1789.7333
$ws.Range("K122").Value = 5369.199900000001
$ws.Range("M122").Value = -2919.199900000001
$ws.Range("H132").Value = 2115.8
$ws.Range("I132").Value = 1999.6666
$ws.Range("K132").Value = 5998.9998
$ws.Range("M132").Value = -3468.9998

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 47070.645
$ws.Range("I122").Value = 56899.434
$ws.Range("K122").Value = 170698.302
$ws.Range("M122").Value = -168248.302
$ws.Range("H123").Value = 47335.7
$ws.Range("J123").Value = 47335.7
$ws.Range("L123").Value = 47335.7
$ws.Range("N123").Value = -57135.7
$ws.Range("H125").Value = 40000
$ws.Range("J125").Value = 40000
$ws.Range("L125").Value = 40000
$ws.Range("N125").Value = -49840
$ws.Range("H126").Value = 4126.3335
$ws.Range("I126").Value = 3171.8096
$ws.Range("K126").Value = 9515.4288
$ws.Range("M126").Value = -7045.4288
$ws.Range("H132").Value = 3299.7778
$ws.Range("I132").Value = 1633.3334
$ws.Range("K132").Value = 4900.0002
$ws.Range("M132").Value = -2370.0002
